$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header row: "Base case scenario" -> "Business-as-usual scenario"
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "Business-as-usual scenario"

# ---------------------------------------------------------------------------
# 2) Insert a new blank spacer row after row 5 (before "2. Increase ...").
#    Formatting is inherited from the row above, matching styles 4/5/5.
# ---------------------------------------------------------------------------
$ws.Rows("6:6").Insert()

# ---------------------------------------------------------------------------
# 3) Update the "Reduce use ..." rows (now rows 8 and 9) with new figures.
# ---------------------------------------------------------------------------
$ws.Range("B8").Value = "74.5% of forage fish to aquaculture "
$ws.Range("C8").Value = "74.5% of forage fish to aquaculture"

$ws.Range("B9").Value = "28.5% of this supply to mariculture"
$ws.Range("C9").Value = "100% of this supply to mariculture"

# ---------------------------------------------------------------------------
# 4) Insert two new rows (10 and 11) for the new "Percent of forage fish ..."
#    line plus its trailing blank spacer row.
# ---------------------------------------------------------------------------
$ws.Rows("10:11").Insert()

$ws.Range("A10").Value = "Percent of forage fish destined for reduction to mariculture:"
$ws.Range("B10").Value = "21.2% of forage fish to mariculture"
$ws.Range("C10").Value = "74.5% of forage fish to mariculture"
$ws.Range("C10").NumberFormat = "0.00%"
$ws.Range("C11").NumberFormat = "0.00%"

# ---------------------------------------------------------------------------
# 5) Update the "3. Reduce the amount of fish ingredients in feed" row
#    (now row 12) with the new 2030/2050 column headers.
# ---------------------------------------------------------------------------
$ws.Range("B12").Value = "2030 FM/FO compositions"
$ws.Range("C12").Value = "2050 FM/FO compositions"

# ---------------------------------------------------------------------------
# 6) Insert a new blank spacer row (13) before the last data row.
# ---------------------------------------------------------------------------
$ws.Rows("13:13").Insert()

# ---------------------------------------------------------------------------
# 7) Update the "4. Reduce the feed conversion rate ..." row (now row 14)
#    with the new 2030/2050 column headers.
# ---------------------------------------------------------------------------
$ws.Range("B14").Value = "2030 feed conversion rates"
$ws.Range("C14").Value = "2050 feed conversion rates"

# ---------------------------------------------------------------------------
# 8) Update the sheet's active-cell selection, matching the saved view.
# ---------------------------------------------------------------------------
$ws.Range("B17").Select()
